$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.977.01'
$ws.Range("E2").Value = '  -3.48%  '

$ws.Range("D3").Value = '''1.715.57'
$ws.Range("E3").Value = '  -3.02%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''308.11'
$ws.Range("E5").Value = '  -6.29%  '

$ws.Range("D7").Value = '''0.4772'
$ws.Range("E7").Value = '  +4.74%  '

$ws.Range("D8").Value = '''0.3478'
$ws.Range("E8").Value = '  -1.37%  '

$ws.Range("D9").Value = '''41.96'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").Value = '''0.07229'
$ws.Range("E10").Value = '  -2.11%  '

$ws.Range("D11").Value = '''1.045'
$ws.Range("E11").Value = '  -4.69%  '

$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("D13").Value = '''19.86'
$ws.Range("E13").Value = '  -4.13%  '

$ws.Range("D14").Value = '''5.838'
$ws.Range("E14").Value = '  -2.81%  '

$ws.Range("D15").Value = '''1.703.50'
$ws.Range("E15").Value = '  -3.87%  '

$ws.Range("D16").Value = '''6.837'
$ws.Range("E16").Value = '  -4.93%  '

$ws.Range("D17").Value = '''86.51'
$ws.Range("E17").Value = '  -6.61%  '

$ws.Range("E18").Value = '  -2.20%  '

$ws.Range("D19").Value = '''0.06378'
$ws.Range("E19").Value = '  -0.88%  '

$ws.Range("D21").Value = '''16.50'
$ws.Range("E21").Value = '  -2.50%  '

$ws.Range("E22").Value = '  -2.78%  '

$ws.Range("D23").Value = '''27.036.52'
$ws.Range("E23").Value = '  -3.36%  '

$ws.Range("D24").Value = '''10.75'
$ws.Range("E24").Value = '  -4.17%  '

$ws.Range("D25").Value = '''2.089'
$ws.Range("E25").Value = '  -2.03%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''151.44'
$ws.Range("E26").Value = '  -5.67%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''19.95'
$ws.Range("E27").Value = '  -0.97%  '

$ws.Range("D28").Value = '''1.909.24'
$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("D29").Value = '''2.083'
$ws.Range("E29").Value = '  -4.56%  '

$ws.Range("D30").Value = '''120.84'
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("E31").Value = '  -4.68%  '

$ws.Range("D32").Value = '''0.09128'
$ws.Range("E32").Value = '  -1.36%  '

$ws.Range("D33").Value = '''3.600'
$ws.Range("E33").Value = '  -1.67%  '

$ws.Range("D34").Value = '''5.334'
$ws.Range("E34").Value = '  -4.75%  '

$ws.Range("D35").Value = '''1.467'
$ws.Range("E35").Value = '  +6.16%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '''0.02178'
$ws.Range("E36").Value = '  -4.61%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.05872'
$ws.Range("E37").Value = '  -4.22%  '

$ws.Range("D38").Value = '''11.00'
$ws.Range("E38").Value = '  -7.14%  '

$ws.Range("D39").Value = '''0.2001'
$ws.Range("E39").Value = '  -4.21%  '

$ws.Range("D40").Value = '''0.6057'
$ws.Range("E40").Value = '  -3.27%  '

$ws.Range("D41").Value = '''4.731'
$ws.Range("E41").Value = '  -3.87%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.085'
$ws.Range("E42").Value = '  -8.23%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''7.450'
$ws.Range("E43").Value = '  -4.97%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''12.70'
$ws.Range("E44").Value = '  -4.00%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '''3.568'
$ws.Range("E45").Value = '  -4.43%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.5630'
$ws.Range("E46").Value = '  -3.91%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''119.03'
$ws.Range("E47").Value = '  -2.91%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.832'
$ws.Range("E48").Value = '  -5.45%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '''1.111'
$ws.Range("E49").Value = '  -1.81%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.06650'
$ws.Range("E50").Value = '  -2.46%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '''1.001'
$ws.Range("E51").Value = '  +0.13%  '

